$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 950; this shifts the existing rows 950..1008
# down to 951..1009 (matching the target dimension A1:T1009).
$ws.Rows(950).Insert()

# Populate the newly inserted row 950 with the new data record.
$ws.Cells.Item(950, 1).Value = 5
$ws.Cells.Item(950, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(950, 3).Value = "Maule"
$ws.Cells.Item(950, 4).Value = 44610
$ws.Cells.Item(950, 5).Value = 7
$ws.Cells.Item(950, 6).Value = "Fruta"
$ws.Cells.Item(950, 7).Value = 100102
$ws.Cells.Item(950, 8).Value = "Cítricos"
$ws.Cells.Item(950, 9).Value = 100102003
$ws.Cells.Item(950, 10).Value = "Limón"
$ws.Cells.Item(950, 11).Value = "Sin especificar"
$ws.Cells.Item(950, 12).Value = "1a plateado"
$ws.Cells.Item(950, 13).Value = 300
$ws.Cells.Item(950, 14).Value = 17000
$ws.Cells.Item(950, 15).Value = 17000
$ws.Cells.Item(950, 16).Value = 17000
$ws.Cells.Item(950, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(950, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(950, 19).Value = 1062
$ws.Cells.Item(950, 20).Value = 16
